$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds numeric-looking strings (e.g. "235.09") that must
# stay as TEXT, matching the original inline-string cells. Pre-format those cells
# as Text so the COM layer does not silently coerce the assigned strings to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "30.327.77"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.869.92"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "235.09"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "0.4708"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").Value = "0.06578"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").Value = "21.65"
$ws.Range("E10").Value = "  -3.27%  "
$ws.Range("D11").Value = "0.08017"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").Value = "97.02"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "1.868.54"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "5.116"
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("D15").Value = "0.6842"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").Value = "269.62"
$ws.Range("E16").Value = "  -3.52%  "
$ws.Range("D17").Value = "30.318.99"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "14.01"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D19").Value = "0.000007619"
$ws.Range("E19").Value = "  +3.82%  "
$ws.Range("D21").Value = "2.113.68"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "5.289"
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("D24").Value = "6.218"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").Value = "9.429"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("D28").Value = "1.948"
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D30").Value = "0.09943"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("D31").Value = "4.365"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("E32").Value = "  -1.40%  "
$ws.Range("D33").Value = "4.071"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "0.04708"
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("D35").Value = "1.140"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "0.7009"
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "0.01869"
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("D39").Value = "2.634"
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("D40").Value = "6.297"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").Value = "71.93"
$ws.Range("E41").Value = "  -6.45%  "
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Value = "0.8421"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").Value = "0.4170"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "102.85"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").Value = "9.185"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("D48").Value = "7.046"
$ws.Range("E48").Value = "  -2.56%  "
$ws.Range("D49").Value = "914.15"
$ws.Range("E49").Value = "  -5.18%  "
$ws.Range("D50").Value = "34.49"
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("E51").Value = "  +0.90%  "

# Restore the default "Normal" style on the Price cells so the Text number format
# applied above does not linger as an explicit cell style (keeps styling identical
# to the original workbook, only the values themselves change).
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
